$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.189.18"
$ws.Range("E2").Value = "  -1.48%  "
# Row 3
$ws.Range("D3").Value = "1.858.09"
$ws.Range("E3").Value = "  -1.13%  "
# Row 4
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.21%  "
# Row 5
$ws.Range("D5").Value = "'232.68"
$ws.Range("E5").Value = "  -2.68%  "
# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.12%  "
# Row 7
$ws.Range("D7").Value = "'0.4740"
$ws.Range("E7").Value = "  -1.69%  "
# Row 8
$ws.Range("D8").Value = "'0.2739"
$ws.Range("E8").Value = "  -3.33%  "
# Row 9
$ws.Range("D9").Value = "'0.06410"
$ws.Range("E9").Value = "  -1.80%  "
# Row 10
$ws.Range("D10").Value = "1.856.02"
$ws.Range("E10").Value = "  -6.32%  "
# Row 11
$ws.Range("D11").Value = "'0.07448"
$ws.Range("E11").Value = "  -0.86%  "
# Row 12
$ws.Range("D12").Value = "'16.25"
$ws.Range("E12").Value = "  -1.84%  "
# Row 13
$ws.Range("D13").Value = "'5.005"
$ws.Range("E13").Value = "  -1.78%  "
# Row 14
$ws.Range("D14").Value = "'85.23"
$ws.Range("E14").Value = "  -4.03%  "
# Row 15
$ws.Range("D15").Value = "'0.6312"
$ws.Range("E15").Value = "  -5.02%  "
# Row 16
$ws.Range("D16").Value = "30.135.59"
$ws.Range("E16").Value = "  -1.54%  "
# Row 17
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.13%  "
# Row 18
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'231.32"
$ws.Range("E18").Value = "  +0.41%  "
# Row 19
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'12.74"
$ws.Range("E19").Value = "  -4.36%  "
# Row 20
$ws.Range("D20").Value = "'0.000007331"
$ws.Range("E20").Value = "  -3.54%  "
# Row 21
$ws.Range("D21").Value = "2.098.42"
$ws.Range("E21").Value = "  -5.31%  "
# Row 22
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.30%  "
# Row 23
$ws.Range("D23").Value = "'5.066"
$ws.Range("E23").Value = "  -4.10%  "
# Row 24
$ws.Range("B24").Value = "BitDAO"
$ws.Range("C24").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D24").Value = "'0.3918"
$ws.Range("E24").Value = "  -7.92%  "
# Row 25
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'5.995"
$ws.Range("E25").Value = "  -2.89%  "
# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.253"
$ws.Range("E26").Value = "  -0.77%  "
# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'166.01"
$ws.Range("E27").Value = "  -1.03%  "
# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.81"
$ws.Range("E28").Value = "  -4.26%  "
# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'1.882"
$ws.Range("E29").Value = "  -3.14%  "
# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.383"
$ws.Range("E30").Value = "  -3.03%  "
# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.1003"
$ws.Range("E31").Value = "  +4.99%  "
# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.120"
$ws.Range("E32").Value = "  -5.12%  "
# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.915"
$ws.Range("E33").Value = "  -3.04%  "
# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04892"
$ws.Range("E34").Value = "  -2.62%  "
# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.141"
$ws.Range("E35").Value = "  -5.75%  "
# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7218"
$ws.Range("E36").Value = "  -3.57%  "
# Row 37
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'0.9998"
$ws.Range("E37").Value = "  -1.02%  "
# Row 38
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.690"
$ws.Range("E38").Value = "  -0.51%  "
# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01887"
$ws.Range("E39").Value = "  +1.97%  "
# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.635"
$ws.Range("E40").Value = "  +0.21%  "
# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8988"
$ws.Range("E41").Value = "  -1.93%  "
# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.952"
$ws.Range("E42").Value = "  -6.42%  "
# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'105.64"
$ws.Range("E43").Value = "  -0.67%  "
# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "  -0.12%  "
# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.561"
$ws.Range("E45").Value = "  -4.41%  "
# Row 46
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4090"
$ws.Range("E46").Value = "  -4.39%  "
# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.035"
$ws.Range("E47").Value = "  -5.18%  "
# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'60.98"
$ws.Range("E48").Value = "  -5.74%  "
# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1195"
$ws.Range("E49").Value = "  -7.04%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.731"
$ws.Range("E50").Value = "  -2.52%  "
# Row 51
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'33.19"
$ws.Range("E51").Value = "  -2.10%  "
